$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "What is the next number in the following sequence?",
        "ques_type": 2,
        "options": [
            "124",
            "139",
            "171",
            "193"
        ],
        "score": "193"
    },
    {
        "title": "Fifteen people working 5 hours per day can make 30 units of a product in 10 days. Assume that all other factors remain constant, and people of the same efficiency are used to make the same products.In how many days can 10 people make 10 units of the product if each of them works 10 hours per day?",
        "ques_type": 2,
        "options": [
            "2.5 days",
            "7.5 days",
            "12 days",
            "26 days"
        ],
        "score": "2.5 days"
    },
    {
        "title": "The table below shows a company\u2019s Manufacturing Cost, Overhead, Total Sales, Profit, and Dividend per Shareholder over 4 years. Assume the relationships among Manufacturing Cost, Overhead, Total Sales, Profit, and Dividend per Shareholder remain the same over the years.What should have been the Dividend per Shareholder in 2017, assuming the number of Shareholders has remained unchanged during the period 2017-2020?",
        "ques_type": 2,
        "options": [
            "$17.45",
            "$19.50",
            "$20.00",
            "$25.00"
        ],
        "score": "$20.00"
    },
    {
        "title": "The figure below depicts a company\u2019s customer analysis based on gender and age group.What percentage of age-groups has more male customers than female customers?",
        "ques_type": 2,
        "options": [
            "37.5%",
            "40%",
            "50%",
            "62.5%"
        ],
        "score": "50%"
    }
]
'@

# Write the long multi-line text into a scratch cell far away from row 1,
# then Cut/Paste it into A1. Writing directly into A1 (or any cell in row 1)
# triggers an automatic row-height recalculation because the text contains
# newlines; performing the write on a scratch row and then moving the cell
# via Cut/Paste avoids that recalculation on row 1.
$scratch = $ws.Range("Z100")
$scratch.Value2 = $questionsText
$scratch.Cut($ws.Range("A1"))

# Remove the now-empty scratch row entirely so it leaves no trace.
$ws.Rows(100).Delete()

# Remove the second row (previously holding the duplicate shared-string cell).
$ws.Rows(2).Delete()

# Clear any direct formatting left on A1 (bold font / border / alignment)
# so the cell reverts to the default style.
$ws.Range("A1").ClearFormats()

Write-Output ("A1 length: " + $ws.Range("A1").Value2.Length)
